# Applies the "Add files via upload" edit to prompts.xlsx:
#  - De-personalise the hyper-personalised prompt (B2) by stripping the
#    {CompanyName}/{feature x}/{feature y}/{objectives} placeholders.
#  - Add three new CSM-feature-recommendation prompts in rows 4-6, numbering
#    them 3, 4, 5 in column A.
#  - Give the newly appended row (row 6) a plain wrap-text style instead of
#    the existing justified Times New Roman style.
#  - Leave the "tailored solutions" prompt (row 3) as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: strip the placeholders from the hyper-personalised prompt -----
$ws.Range("B2").Value = 'You''re an experienced and customer-focused customer success manager with exceptional communication skills. For many years, you have handled client correspondence and emails for various projects, ensuring clarity and professionalism in all your interactions. Your task is to draft an email to a client serving as an alert based on the recommendation model, advising the company on the underutilisation of features from your company''s product. This email has to be extremely hyper-personalised. Specifically, recommend incorporating features to enhance efficiency, productivity, customer satisfaction and reach their objectives. Please ensure that the tone of the email is friendly yet professional, providing all necessary information clearly and concisely. Include any relevant attachments or links the client may need for reference or further action.'

# --- Row 3 (prompt_id 2) is unchanged - "tailored solutions" prompt stays -

# --- Row 4 (prompt_id 3): new "onboarded but underused feature" prompt ----
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'You, as a Customer success Manager, are responsible for ensuring clients make the most of your product. Draft an email to a client who has been onboarded but hasn''t fully utilized a key feature that could significantly impact their productivity. Start the email by expressing your commitment to their success and referencing specific data points or observations from their usage patterns. Describe the underused feature, detail its direct benefits, and draw parallels with similar clients who have seen improvements by adopting this feature. Offer to set up a personalized training session specifically for their team. Conclude the email by suggesting a short call to discuss their current experience and how this feature fits into their workflow'

# --- Row 5 (prompt_id 4): new "noticed underutilised feature" prompt ------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'You are a Customer Success Manager who has noticed that a specific client could improve their efficiency by using an underutilized feature of your product. Compose an email where you introduce this feature in detail. Begin by acknowledging their current use of your product and complimenting them on their successes. Explain the feature, focusing on how it integrates with their existing workflow and the specific benefits it offers. Include data-driven insights and show some performance metrics to create a strong impact. Offer to provide a personalized demo or training session to help them get started. Close the email by suggesting a follow-up meeting to discuss implementation.'

# --- Row 6 (prompt_id 5): new "proactive CSM" prompt, plain wrap style ----
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'As a proactive Customer Success Manager, you''ve spotted a chance to boost a client''s results through a feature recommendation. Compose an email to introduce this feature comprehensively. Start with a warm acknowledgment of how well they''ve adapted to your core offerings and their innovative use of your product. Describe the neglected feature, focusing on its synergy with their current workflow and how it can solve specific problems or speed up processes. Share empirical evidence or analytics that illustrate the potential impact of adopting this feature. Offer a personalized walkthrough or interactive demo tailored to their business needs. Wrap up the email by proposing a meeting to discuss how this feature can be implemented to maximize their productivity and outcomes'
$ws.Range("B6").WrapText = $true

# --- Row heights to match the re-flowed (shorter/longer) text content -----
$ws.Rows.Item(2).RowHeight = 234
$ws.Rows.Item(4).RowHeight = 202.8
$ws.Rows.Item(5).RowHeight = 187.2
$ws.Rows.Item(6).RowHeight = 187.2

# --- Selection / scroll position, matching the saved view -----------------
$ws.Application.GoTo($ws.Range("A5"))
$ws.Range("B6").Select()
